$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Create github account...
$ws.Range('C2').Value = 1
$ws.Range('D2').Value = "Completed"
$ws.Range('E2').Value = [DateTime]"2016-06-07"
$ws.Range('F2').Value = [DateTime]"2016-06-09"

# Row 3 - Play around with drone...
$ws.Range('C3').Value = 0.5
$ws.Range('D3').Value = "In progress"
$ws.Range('E3').Value = [DateTime]"2016-06-08"

# Row 4 - Experiment with different recording modes...
$ws.Range('C4').Value = 1.5
$ws.Range('D4').Value = "In progress"
$ws.Range('E4').Value = [DateTime]"2016-06-08"

# Row 5 - Figure out battery life between charge...
$ws.Range('C5').Value = 0.5
$ws.Range('D5').Value = "Completed"
$ws.Range('E5').Value = [DateTime]"2016-06-08"
$ws.Range('F5').Value = [DateTime]"2016-06-08"
$ws.Range('H5').Value = "Approximately 15 mins. Does not allow for takeoff after dropping below ~30%"

# Row 4 Notes, filled in after row 5 (matches the shared-string append order)
$ws.Range('H4').Value = "Drone does not follow vertically, continues at about its same height from a distance while on ""leash"" mode."

# Move active selection to H3 to mirror the author's last-saved cursor position
$ws.Range("H3").Select()
